# project.xlsx - mark the "分析功能二" (analysis feature 2) row as complete
# and make the "是否完成" (completed?) column a Yes/No dropdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F7 ("是否完成" for 分析功能二) flips from "否" (No) to "是" (Yes)
$ws.Range("F7").Value = "是"

# Give the whole "是否完成" column (F4:F7) a Yes/No list validation
$ws.Range("F4:F7").Validation.Add(3, 1, 1, """是,否""")

# Leave the selection on the cell that was last touched
$ws.Range("F6").Select() | Out-Null
